$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 0.8459207920792079
$ws.Range("C8").Value = 0.04455585182996977
$ws.Range("D8").Value = 0.8059709104855858
$ws.Range("E8").Value = 0.04603503786206437
$ws.Range("F8").Value = 0.9152418300653596
$ws.Range("G8").Value = 0.07165552846278239
$ws.Range("H8").Value = 0.8553178420801515
$ws.Range("I8").Value = 0.04329404046302684
$ws.Range("J8").Value = 0.767881188118812
$ws.Range("K8").Value = 0.04230384303784946
$ws.Range("L8").Value = 0.7440600604448825
$ws.Range("M8").Value = 0.06677733342194654
$ws.Range("N8").Value = 0.8425620915032679
$ws.Range("O8").Value = 0.1258679662470351
$ws.Range("P8").Value = 0.7803789327126243
$ws.Range("Q8").Value = 0.0529641708465006
$ws.Range("B9").Value = 0.9253069306930695
$ws.Range("C9").Value = 0.02696885764030971
$ws.Range("D9").Value = 0.880477462780817
$ws.Range("E9").Value = 0.03885462910133234
$ws.Range("F9").Value = 0.9868496732026144
$ws.Range("G9").Value = 0.02550240699645172
$ws.Range("H9").Value = 0.9300022840041442
$ws.Range("I9").Value = 0.02440057492884872
$ws.Range("J9").Value = 0.8823762376237624
$ws.Range("K9").Value = 0.03841772987091333
$ws.Range("L9").Value = 0.8298394374351592
$ws.Range("M9").Value = 0.0417845885789845
$ws.Range("N9").Value = 0.9644183006535948
$ws.Range("O9").Value = 0.03396886397771902
$ws.Range("P9").Value = 0.891670847266695
$ws.Range("Q9").Value = 0.03423675554348856
